$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove all borders from the whole used area (header + new body rows) ---
$ws.Range("A1:H5").Borders.LineStyle = 0

# --- Header row: drop the old "vertical top" alignment, keep left-align on column C ---
$ws.Range("A1:H1").VerticalAlignment = -4107
$ws.Range("C1").HorizontalAlignment = -4131
$ws.Range("C1:C5").VerticalAlignment = -4107

# --- Rewrite row 2 (previously the single "Credential Issuance" test case) ---
$ws.Range("A2").Value = "Resident Service_Best two Fingers_01"
$ws.Range("B2").Value = "Resident Services"
$ws.Range("C2").Value = "Best two fingers"
$ws.Range("D2").Value = "Functional"
$ws.Range("E2").Value = "Verify QR code having the best two finger"
$ws.Range("F2").Value = "QR code should have best two fingers"
$ws.Range("G2").Value = "N"

# --- Add new row 3 ---
$ws.Range("A3").Value = "Resident Service_Best two Fingers_02"
$ws.Range("B3").Value = "Resident Services"
$ws.Range("C3").Value = "Best two fingers"
$ws.Range("D3").Value = "Functional"
$ws.Range("E3").Value = "Verify QR code having the best fingers with having only specific fingers in policy"
$ws.Range("F3").Value = "QRcode should have best two fingers based on fingers mentioned in policy"
$ws.Range("G3").Value = "N"

# --- Add new row 4 ---
$ws.Range("A4").Value = "Resident Service_Best two Fingers_03"
$ws.Range("B4").Value = "Resident Services"
$ws.Range("C4").Value = "Best two fingers"
$ws.Range("D4").Value = "Functional"
$ws.Range("E4").Value = "Verify QR code having the best fingers by not giving specific fingers in policy"
$ws.Range("F4").Value = "QRcode should have best two fingers by comparing score of all the fingers"
$ws.Range("G4").Value = "N"

# --- Add new row 5 ---
$ws.Range("A5").Value = "Resident Service_Best two Fingers_04"
$ws.Range("B5").Value = "Resident Services"
$ws.Range("C5").Value = "Best two fingers"
$ws.Range("D5").Value = "Functional"
$ws.Range("E5").Value = "Verify rank for best fingers given based on score"
$ws.Range("F5").Value = "Rank for best fingers should be based on score and rank one finger should be first and rank two finger should be second"
$ws.Range("G5").Value = "N"

# --- Formatting: wrap text + explicit black font colour on the new body data,
#     matching every column except the (unstyled) Type/Reviewed columns D & H ---
$ws.Range("A2:C5").Font.Color = 0
$ws.Range("A2:C5").WrapText = $true
$ws.Range("E2:G5").Font.Color = 0
$ws.Range("E2:G5").WrapText = $true

# Column widths / heights
$ws.Columns.Item(5).ColumnWidth = 30.5
$ws.Rows.Item(2).RowHeight = 43.5
$ws.Rows.Item(3).RowHeight = 46.5
$ws.Rows.Item(4).RowHeight = 55.5
$ws.Rows.Item(5).RowHeight = 63.75

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection (matches the saved cursor position in the workbook)
[void]$ws.Range("H4").Select()
